# El-Maven v0.11.0 added a new `adductName` output column. Insert it into
# every worksheet (as a non-sample column) so TraceBase stops treating it
# as sample data.

$wb = $excel.ActiveWorkbook

# "Original" sheet: adductName goes right before isotopeLabel (column H).
$ws = $wb.Worksheets.Item("Original")
$ws.Columns("H:H").Insert()
$ws.Range("H1").Value = "adductName"
$ws.Range("H2:H12").Value = "[M-H]-"

# "Corrected" sheet: adductName goes right after C_Label (column C).
$ws = $wb.Worksheets.Item("Corrected")
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "adductName"
$ws.Range("C2:C12").Value = "[M-H]-"

# "Normalized" sheet: same layout as Corrected.
$ws = $wb.Worksheets.Item("Normalized")
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "adductName"
$ws.Range("C2:C12").Value = "[M-H]-"

# "PoolAfterDF" sheet: adductName goes right after Compound (column B).
$ws = $wb.Worksheets.Item("PoolAfterDF")
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "adductName"
$ws.Range("B2:B3").Value = "[M-H]-"
